$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------
# Sheet2 formatting: grab the two label/content style pairs (currently
# sitting on Sheet1 rows 11/12 - the "Method C25 convertStringToC25"
# callout) BEFORE Sheet1 gets wiped, since those are the only cells in
# the workbook carrying styles 6/7 (wrapped content box).
# ---------------------------------------------------------------------
$ws1.Range("C11:D11").Copy() | Out-Null
$ws2.Range("B3:C3").PasteSpecial(-4122) | Out-Null
$ws1.Range("C12:D12").Copy() | Out-Null
$ws2.Range("B4:C4").PasteSpecial(-4122) | Out-Null
$ws1.Range("C11:D11").Copy() | Out-Null
$ws2.Range("B6:C6").PasteSpecial(-4122) | Out-Null
$ws1.Range("C12:D12").Copy() | Out-Null
$ws2.Range("B7:C7").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# Sheet1: the header band (rows 3-5) stays untouched. Everything from
# row 6 down (the sample mapping rows + the two method-signature
# "callout" blocks) gets wiped and rebuilt.
# ---------------------------------------------------------------------
$ws1.Range("C6:H22").UnMerge() | Out-Null
$ws1.Range("C6:H22").Clear() | Out-Null

# drop the stray custom-height rows left behind by the old callout blocks
$ws1.Rows.Item(12).EntireRow.AutoFit() | Out-Null
$ws1.Rows.Item(15).EntireRow.AutoFit() | Out-Null
$ws1.Rows.Item(19).EntireRow.AutoFit() | Out-Null

# -- data rows --------------------------------------------------------
$ws1.Range("C6").Value = "ChPolicyEntity"
$ws1.Range("D6").Value = "HomePolicyQuoteInqRqType"
$ws1.Range("E6").Value = "policyNumber"
$ws1.Range("F6").Value = "persPolicy.policyNumber"
$ws1.Range("G6").Value = "convertStringToC25"

$ws1.Range("C7").Value = "ChPolicyEntity"
$ws1.Range("D7").Value = "HomePolicyQuoteInqRqType"
$ws1.Range("E7").Value = "version"
$ws1.Range("F7").Value = "persPolicy.policyVersion.stringValue"

$ws1.Range("C10").Value = "ChPolicyEntity"
$ws1.Range("D10").Value = "HomePolicyQuoteInqRqType"
$ws1.Range("E10").Value = "contractTerm.effective"
$ws1.Range("F10").Value = "persPolicy.contractTerm.effectiveDt.stringValue"

$ws1.Range("C11").Value = "ChPolicyEntity"
$ws1.Range("D11").Value = "HomePolicyQuoteInqRqType"
$ws1.Range("E11").Value = "contractTerm.expiration"
$ws1.Range("F11").Value = "persPolicy.contractTerm.expirationDt.stringValue"

$ws1.Range("C12").Value = "ChPolicyEntity"
$ws1.Range("D12").Value = "HomePolicyQuoteInqRqType"
$ws1.Range("E12").Value = "producerCd"
$ws1.Range("F12").Value = "producerArray[0].producerInfo.contractNumber.stringValue"

$ws1.Range("C13").Value = "ChPolicyEntity"
$ws1.Range("D13").Value = "HomePolicyQuoteInqRqType"
$ws1.Range("E13").Value = "producerCd"
$ws1.Range("F13").Value = "producerArray[1].producerInfo.contractNumber.stringValue"

# two trailing blank spacer rows (kept from the old callout-block rows,
# now empty except for their custom height)
$ws1.Rows.Item(14).RowHeight = 80.25
$ws1.Rows.Item(17).RowHeight = 24

$ws1.Range("C8").Select() | Out-Null

Write-Host "sheet1 done"

# ---------------------------------------------------------------------
# Sheet2: used to be completely empty; it now hosts the two
# "convertStringToC25" method-signature / code callouts that used to
# live at the bottom of Sheet1 (rows 11/12 and 14/15), re-flowed into a
# two-column (B/C) layout.
# ---------------------------------------------------------------------
$ws2.Columns.Item(3).ColumnWidth = 49.42578125

$ws2.Range("B3").Value = "Method C25 convertStringToC25(String source, C25 destination)"
$ws2.Range("B3:C3").Merge() | Out-Null
$ws2.Rows.Item(3).RowHeight = 15

$ws2.Range("B4").Value = "if (destination==null) {`n   destination = (C25)XmlBeanFactory.newInstance(C25.class);`n}`ndestination.setId(source);`nreturn destination;"
$ws2.Range("B4:C4").Merge() | Out-Null
$ws2.Rows.Item(4).RowHeight = 99.75

$ws2.Range("B6").Value = "Method String convertStringToC25(C25 source, String destination)"
$ws2.Range("B6:C6").Merge() | Out-Null
$ws2.Rows.Item(6).RowHeight = 15

$ws2.Range("B7").Value = "`nreturn source.getId();"
$ws2.Range("B7:C7").Merge() | Out-Null
$ws2.Rows.Item(7).RowHeight = 24.75

$ws2.Range("C10").Select() | Out-Null

Write-Host "sheet2 done"
